$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 5.3
$ws.Range("G2").Value = 9
$ws.Range("H2").Value = 1.5
$ws.Range("I2").Value = 1.6
$ws.Range("J2").Value = 3.95
$ws.Range("K2").Value = 5.3
$ws.Range("L2").Value = 1.37
$ws.Range("M2").Value = 1.06
$ws.Range("N2").Value = 3.25
$ws.Range("O2").Value = 1.3
$ws.Range("P2").Value = 1.86
$ws.Range("Q2").Value = 1.79
$ws.Range("R2").Value = 1.35
$ws.Range("S2").Value = 2.92
$ws.Range("T2").Value = 1.98
$ws.Range("U2").Value = 1.78
$ws.Range("V2").Value = 2.66
$ws.Range("W2").Value = 1.14
$ws.Range("X2").Value = 19
$ws.Range("Z2").Value = 10.5
$ws.Range("AA2").Value = 16.5
$ws.Range("AB2").Value = 26
$ws.Range("AC2").Value = 11
$ws.Range("AD2").Value = 12
$ws.Range("AE2").Value = 21
$ws.Range("AF2").Value = 75
$ws.Range("AG2").Value = 34
$ws.Range("AK2").Value = 1000
$ws.Range("AO2").Value = 10.5
$ws.Range("J3").Value = 1.01
$ws.Range("N3").Value = 1.02
$ws.Range("P3").Value = 1.25
$ws.Range("R3").Value = 1.25
$ws.Range("AC4").Value = 9.800000000000001
$ws.Range("F5").Value = 1.05
$ws.Range("G5").Value = 980
$ws.Range("H5").Value = 1.05
$ws.Range("J5").Value = 1.06
$ws.Range("H6").Value = 3
$ws.Range("O6").Value = 1.32
$ws.Range("P6").Value = 1.75
$ws.Range("I7").Value = 1.84
$ws.Range("N7").Value = 1.63
$ws.Range("P7").Value = 1.63
$ws.Range("Q7").Value = 1.86
$ws.Range("S7").Value = 1.86
$ws.Range("V7").Value = 2.18
$ws.Range("F8").Value = 1.81
$ws.Range("G8").Value = 2.1
$ws.Range("H8").Value = 4.6
$ws.Range("I8").Value = 6.2
$ws.Range("J8").Value = 3.05
$ws.Range("K8").Value = 3.95
$ws.Range("L8").Value = 1.01
$ws.Range("N8").Value = 2.62
$ws.Range("O8").Value = 1.45
$ws.Range("T8").Value = 2.04
$ws.Range("V8").Value = 1.19
$ws.Range("W8").Value = 1.92
$ws.Range("Y8").Value = 15
$ws.Range("AB8").Value = 7.4
$ws.Range("AC8").Value = 8.199999999999999
$ws.Range("AF8").Value = 12
$ws.Range("AG8").Value = 12.5
$ws.Range("AJ8").Value = 27
$ws.Range("AK8").Value = 30
$ws.Range("AN8").Value = 26
$ws.Range("H9").Value = 2.74
$ws.Range("K9").Value = 3.75
$ws.Range("O9").Value = 1.32
$ws.Range("P9").Value = 1.87
$ws.Range("S9").Value = 3.2
$ws.Range("J10").Value = 3.5
$ws.Range("F11").Value = 1.62
$ws.Range("I11").Value = 5.8
$ws.Range("S11").Value = 2.24
$ws.Range("V11").Value = 1.21
$ws.Range("AH11").Value = 1000
$ws.Range("AN11").Value = 6.2
$ws.Range("F12").Value = 2.7
$ws.Range("H12").Value = 2.6
$ws.Range("J12").Value = 3.65
$ws.Range("K12").Value = 3.8
$ws.Range("N12").Value = 4.5
$ws.Range("Q12").Value = 1.73
$ws.Range("F13").Value = 3.45
$ws.Range("J13").Value = 3.85
$ws.Range("L13").Value = 1.01
$ws.Range("O13").Value = 1.19
$ws.Range("P13").Value = 2.5
$ws.Range("S13").Value = 2.4
$ws.Range("T13").Value = 1.55
$ws.Range("Y13").Value = 14.5
$ws.Range("AC13").Value = 9.6
$ws.Range("AD13").Value = 11.5
$ws.Range("AH13").Value = 1000
$ws.Range("AK13").Value = 44
$ws.Range("AN13").Value = 29
$ws.Range("F14").Value = 6.2
$ws.Range("G14").Value = 6.8
$ws.Range("H14").Value = 1.59
$ws.Range("I14").Value = 1.65
$ws.Range("J14").Value = 4.2
$ws.Range("K14").Value = 4.6
$ws.Range("N14").Value = 5.1
$ws.Range("P14").Value = 2.44
$ws.Range("T14").Value = 1.72
$ws.Range("U14").Value = 2.24
$ws.Range("V14").Value = 2.52
$ws.Range("W14").Value = 1.17
$ws.Range("AA14").Value = 16
$ws.Range("AB14").Value = 28
$ws.Range("AD14").Value = 10.5
$ws.Range("AG14").Value = 23
$ws.Range("AH14").Value = 18.5
$ws.Range("AL14").Value = 70
$ws.Range("AO14").Value = 6.6
$ws.Range("F15").Value = 1.05
$ws.Range("H15").Value = 1.05
$ws.Range("J15").Value = 1.05
